{"js": "// The \"STATUS: \" run (bold label right before \"IN PROGRESS\" in the\n// Feature 1 summary line) is removed, and the document's \"_GoBack\"\n// bookmark (Word's auto-maintained \"last edit position\" marker) is\n// moved from its old spot \u2014 right before the \"(stav)\" parenthetical\n// near the end of the document \u2014 to the spot the deleted run used to\n// occupy.\n//\n// Net effect: delete the old _GoBack bookmark pair where it sits today,\n// then replace the \"STATUS: \" run with a fresh (empty) _GoBack bookmark\n// pair at that location.\n\n// 1) Drop the bookmark from its current location first so only one\n//    \"_GoBack\" exists while we re-insert it (Bookmarks are unique by\n//    name and Document.deleteBookmark() always targets the existing one).\ncontext.document.deleteBookmark(\"_GoBack\");\n\n// 2) Locate the run to remove. It is a unique, exact piece of text in\n//    the document (\"STATUS: \", bold, right before \"IN PROGRESS\").\nconst hits = context.document.body.search(\"STATUS: \", { matchCase: true });\nawait context.sync();\n\nif (hits.items.length !== 1) {\n  throw new Error(`expected exactly one \"STATUS: \" match, found ${hits.items.length}`);\n}\nconst target = hits.items[0];\n\n// 3) Insert the bookmark exactly where \"STATUS: \" starts, then delete\n//    the text \u2014 leaving just <w:bookmarkStart/><w:bookmarkEnd/> behind,\n//    matching the diff.\ntarget.insertBookmark(\"_GoBack\");\ntarget.insertText(\"\", \"Replace\");\n\nawait context.sync();\n", "ps1": "# The \"STATUS: \" run (bold label right before \"IN PROGRESS\" in the\n# Feature 1 summary line) is removed, and the document's \"_GoBack\"\n# bookmark (Word's auto-maintained \"last edit position\" marker) moves\n# from its old spot - right before the \"(stav)\" parenthetical near the\n# end of the document - to the spot the deleted run used to occupy.\n\n$d = $word.ActiveDocument\n\n# Locate the unique \"STATUS: \" run via Find (bold label preceding the\n# \"IN PROGRESS\" status text in the \"Feature 1\" summary paragraph).\n$r = $d.Content\n$found = $r.Find.Execute(\"STATUS: \", $true)\nif (-not $found) {\n    throw \"Could not find the 'STATUS: ' run to remove\"\n}\n\n# Delete the found text, collapsing the range to the insertion point it\n# used to start at.\n$r.Text = \"\"\n\n# Re-seat the (document-unique) \"_GoBack\" bookmark onto that now-empty\n# range. Bookmarks.Add() with an existing name moves the bookmark\n# rather than creating a duplicate, so this both plants the new\n# <w:bookmarkStart/><w:bookmarkEnd/> pair here and removes it from its\n# old location next to \"(stav)\" in one step.\n$d.Bookmarks.Add(\"_GoBack\", $r) | Out-Null\n"}
